# Kubernetes Session notes updated
#
# 1. Delete the slide "Kubernetes Cluster Architecture" / "Refer to draw.io"
#    (the 9th slide) - this also removes its associated notes page.
# 2. Refresh the cached "today" date field text (21-12-2024 -> 21-09-2025)
#    on the Notes Master and every Slide Layout's Date placeholder.

$p = $ppt.ActivePresentation

# --- 1. Remove slide 9 ("Kubernetes Cluster Architecture") ---
$p.Slides.Item(9).Delete()

# --- 2. Update cached date placeholder text across layouts + notes master ---
$oldDate = "21-12-2024"
$newDate = "21-09-2025"

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lyt = $layouts.Item($li)
    for ($i = 1; $i -le $lyt.Shapes.Count; $i++) {
        $sh = $lyt.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
